$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto Price (D) / Volume(1h) (E) cells.
# A leading apostrophe forces Excel to store the value as literal Text
# (matching the workbooks existing inlineStr/text cells) instead of
# re-parsing it as a Number; the Style is then reset to Normal so the
# quote-prefix formatting does not linger on the cell.
$ws.Range('D2').Value = "'27.131.56"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.08%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.900.22"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.05%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.37%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'306.97"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.33%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.29%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.5231"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.51%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3805"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.83%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.07289"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.61%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'21.35"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.90%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.9029"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.56%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.08205"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.60%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'95.47"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.80%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'1.838.32"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.59%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'5.353"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.72%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.31%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.000008670"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.63%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'14.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.13%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.28%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'27.169.27"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.10%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.44%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'2.087.74"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.21%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.04%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'6.446"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.40%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +2.32%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.323"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +2.03%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'18.29"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.08%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.99%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'115.67"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.85%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'4.818"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = "'4.900"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.44%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -0.43%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.05037"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.17%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.7909"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -2.64%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.220"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.16%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.959"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.69%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'3.362"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.80%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'2.643"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +2.33%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.5724"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.26%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.01989"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.51%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.081"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.78%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'9.083"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.65%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'6.611"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.80%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'116.32"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.44%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.1516"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.35%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.4891"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.22%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.30%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'10.18"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.34%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +1.41%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'38.50"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +2.85%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'63.94"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.66%  "
$ws.Range('E51').Style = 'Normal'
